# Refresh the cryptocurrency Price (D) and Volume(1h) (E) columns
# with the latest scraped figures (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @{ D = <new price text, or $null if unchanged>; E = <new volume text> }
$updates = @{
    2 = @{ D = '43.458.04'; E = '  -6.33%  ' }
    3 = @{ D = '2.524.24'; E = '  -3.36%  ' }
    4 = @{ D = $null; E = '  -0.06%  ' }
    5 = @{ D = '297.58'; E = '  -3.43%  ' }
    6 = @{ D = '94.75'; E = '  -5.53%  ' }
    7 = @{ D = '0.575'; E = '  -4.81%  ' }
    8 = @{ D = $null; E = '  +0.13%  ' }
    9 = @{ D = '0.553'; E = '  -4.86%  ' }
    10 = @{ D = '36.47'; E = '  -7.76%  ' }
    11 = @{ D = '0.0805'; E = '  -4.87%  ' }
    12 = @{ D = '7.70'; E = '  -5.73%  ' }
    13 = @{ D = $null; E = '  +0.93%  ' }
    14 = @{ D = '2.915.79'; E = '  -3.16%  ' }
    15 = @{ D = '2.531.43'; E = '  -3.12%  ' }
    16 = @{ D = '0.875'; E = '  -5.04%  ' }
    17 = @{ D = '14.11'; E = '  -5.60%  ' }
    18 = @{ D = '43.446.05'; E = '  -6.63%  ' }
    19 = @{ D = '0.0₃0968'; E = '  -4.49%  ' }
    20 = @{ D = '6.56'; E = '  -2.70%  ' }
    21 = @{ D = '12.34'; E = '  -4.91%  ' }
    22 = @{ D = '72.59'; E = '  +1.54%  ' }
    23 = @{ D = '261.34'; E = '  -4.42%  ' }
    24 = @{ D = '2.91'; E = '  -4.31%  ' }
    25 = @{ D = '2.16'; E = '  -0.28%  ' }
    26 = @{ D = '28.82'; E = '  -1.08%  ' }
    27 = @{ D = '0.999'; E = '  -0.22%  ' }
    28 = @{ D = '10.07'; E = '  -5.17%  ' }
    29 = @{ D = '2.24'; E = '  +0.50%  ' }
    30 = @{ D = '37.32'; E = '  -4.30%  ' }
    31 = @{ D = '6.06'; E = '  -4.57%  ' }
    32 = @{ D = '3.47'; E = '  -4.99%  ' }
    33 = @{ D = '150.78'; E = '  -0.27%  ' }
    34 = @{ D = '2.77'; E = '  -3.08%  ' }
    35 = @{ D = '2.15'; E = '  -4.06%  ' }
    36 = @{ D = '0.0799'; E = '  -4.90%  ' }
    37 = @{ D = '0.115'; E = '  -5.34%  ' }
    38 = @{ D = '0.119'; E = '  -3.66%  ' }
    39 = @{ D = '23.58'; E = '  +0.77%  ' }
    40 = @{ D = '16.21'; E = '  +1.86%  ' }
    41 = @{ D = '3.50'; E = '  -4.26%  ' }
    42 = @{ D = '0.0310'; E = '  -6.63%  ' }
    43 = @{ D = '3.82'; E = '  -6.75%  ' }
    44 = @{ D = '2.018.60'; E = '  -5.04%  ' }
    45 = @{ D = $null; E = '  +0.02%  ' }
    46 = @{ D = '86.38'; E = '  -7.83%  ' }
    47 = @{ D = '1.65'; E = '  +5.52%  ' }
    48 = @{ D = '8.92'; E = '  -6.22%  ' }
    49 = @{ D = '2.777.65'; E = '  -3.13%  ' }
    50 = @{ D = '103.23'; E = '  -5.52%  ' }
    51 = @{ D = '0.188'; E = '  -6.19%  ' }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($null -ne $vals.D) {
        # Prefix with a literal apostrophe (like typing ' in Excel) so that
        # numeric-looking price strings (e.g. "297.58") are stored as text,
        # matching the inline-string cells in the source sheet.
        $ws.Range("D$row").Value = "'" + $vals.D
    }
    if ($null -ne $vals.E) {
        $ws.Range("E$row").Value = $vals.E
    }
}
